# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Swap the Lugo / Almeria rows (names swap, "Casos activos" values swap,
# other columns stay identical so nothing else visibly changes).
$ws.Range("A47").Value = "Almeria"
$ws.Range("C47").Value = 72

$ws.Range("A48").Value = "Lugo"
$ws.Range("C48").Value = 5

# Update the "last updated" timestamp in A1.
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 03:16"
